$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 123
$ws.Range("B3").Value = "satuduatiga"

$range = $ws.Range("A3:B3")
$range.HorizontalAlignment = -4131

$ws.Range("A3:B3").Select()
